$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 4: bonus mission 3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Zoldados de asalto"
$ws.Range("C4").Value = "Estrambotiko"
$ws.Range("D4").Value = "Armadura"

$ws.Hyperlinks.Add($ws.Range("D4"), "Imagenes\Recompensa-Mision-3.png") | Out-Null
$ws.Range("D4").Style = $ws.Range("D3").Style

# Rows 5-21: remaining numbering only (column A)
for ($i = 5; $i -le 21; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

$ws.Range("D11").Select() | Out-Null
